# The sentence "Below are the cross-validation results for both the
# algorithms." was originally split across three separate runs
# (with identical run-level formatting: Comic Sans MS, szCs=22,
# lang=en-US). Re-typing/replacing the full sentence via Find &
# Replace consolidates it back into a single run while keeping the
# shared formatting.

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "Below are the cross-validation results for both the algorithms.",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Below are the cross-validation results for both the algorithms.",
    2
)
